$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update B2:B7 values (+6 offset)
$ws.Range("B2").Value = 19
$ws.Range("B3").Value = 20
$ws.Range("B4").Value = 21
$ws.Range("B5").Value = 22
$ws.Range("B6").Value = 23
$ws.Range("B7").Value = 24

# Update the active selection to B8
$ws.Range("B8").Select()
